$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("2025-05-23 08:42:06", "Parclose", "Sortie - Demande 20250523_082712", 2, 9, 7),
    @("2025-05-23 08:42:06", "Vis 6x50mm", "Sortie - Demande 20250523_082712", 8, 1234, 1226),
    @("2025-05-23 08:42:06", "Tournevis cruciforme", "Sortie - Demande 20250523_082712", 1, 30, 29),
    @("2025-05-23 08:42:06", "Marteau 500g", "Sortie - Demande 20250523_082712", 1, 30, 29)
)

$startRow = 19
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
}
